# Update imputed values in columns D and E for the RandomForest result sheet.
# (Re-run of the imputation algorithm produced slightly different values.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.43140000000001
$ws.Range("D9").Value = -7.519200000000001
$ws.Range("E12").Value = 17.86690000000002
$ws.Range("D13").Value = -8.402800000000001
$ws.Range("E14").Value = 16.78220000000001
$ws.Range("D16").Value = -8.711800000000006
$ws.Range("D18").Value = -8.376300000000001
$ws.Range("E19").Value = 16.3781
$ws.Range("D20").Value = -7.217599999999999
$ws.Range("D26").Value = -8.057199999999993
$ws.Range("E26").Value = 16.39359999999999
$ws.Range("D27").Value = -8.417300000000003
$ws.Range("E27").Value = 16.6404
$ws.Range("D29").Value = -8.017199999999999
$ws.Range("E29").Value = 16.49889999999999
$ws.Range("D35").Value = -7.597100000000001
$ws.Range("D36").Value = -7.648700000000006
$ws.Range("E37").Value = 16.70520000000002
$ws.Range("E38").Value = 16.66850000000001
$ws.Range("D45").Value = -7.509799999999996
$ws.Range("E47").Value = 16.59379999999999
$ws.Range("E51").Value = 17.323
$ws.Range("E52").Value = 17.21780000000001
$ws.Range("D55").Value = -8.975100000000001
$ws.Range("E55").Value = 16.62229999999999
$ws.Range("D57").Value = -8.511500000000002
$ws.Range("D69").Value = -7.880899999999994
$ws.Range("E69").Value = 16.6247
$ws.Range("E70").Value = 17.82250000000002
$ws.Range("D76").Value = -7.626199999999998
$ws.Range("E76").Value = 16.53379999999999
$ws.Range("D78").Value = -7.634700000000002
$ws.Range("E81").Value = 16.3858
$ws.Range("D82").Value = -8.339899999999989
$ws.Range("D83").Value = -9.177900000000001
$ws.Range("E83").Value = 16.5639
$ws.Range("D93").Value = -6.522600000000001
$ws.Range("E94").Value = 18.53460000000003
$ws.Range("D97").Value = -8.207200000000002
$ws.Range("E100").Value = 16.47389999999999
$ws.Range("E102").Value = 16.91289999999999
